$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (LEGO 31173 "Animaux sauvages : le toucan tropical"); all rows below shift up.
$ws.Rows.Item(13).Delete()
